$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '63.652.89'
Set-TextValue $ws.Range('E2') '  +2.73%  '

Set-TextValue $ws.Range('D3') '2.476.34'
Set-TextValue $ws.Range('E3') '  +2.41%  '

Set-TextValue $ws.Range('E4') '  +0.11%  '

Set-TextValue $ws.Range('D5') '576.43'
Set-TextValue $ws.Range('E5') '  +2.45%  '

Set-TextValue $ws.Range('D6') '149.04'
Set-TextValue $ws.Range('E6') '  +4.20%  '

Set-TextValue $ws.Range('E7') '  -0.19%  '

Set-TextValue $ws.Range('E8') '  +1.83%  '

Set-TextValue $ws.Range('E9') '  +5.11%  '

Set-TextValue $ws.Range('E10') '  +0.70%  '

Set-TextValue $ws.Range('E11') '  +2.93%  '

Set-TextValue $ws.Range('E12') '  +3.85%  '

Set-TextValue $ws.Range('D13') '27.50'
Set-TextValue $ws.Range('E13') '  +5.73%  '

Set-TextValue $ws.Range('E14') '  +7.12%  '

Set-TextValue $ws.Range('D15') '2.925.90'
Set-TextValue $ws.Range('E15') '  +2.38%  '

Set-TextValue $ws.Range('D16') '63.526.18'
Set-TextValue $ws.Range('E16') '  +2.68%  '

Set-TextValue $ws.Range('D17') '2.493.24'
Set-TextValue $ws.Range('E17') '  +3.26%  '

Set-TextValue $ws.Range('D18') '11.59'
Set-TextValue $ws.Range('E18') '  +1.95%  '

Set-TextValue $ws.Range('E19') '  +6.81%  '

Set-TextValue $ws.Range('B20') 'Polkadot'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D20') '4.25'
Set-TextValue $ws.Range('E20') '  +2.78%  '

Set-TextValue $ws.Range('B21') 'BitcoinCash'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D21') '329.12'
Set-TextValue $ws.Range('E21') '  +1.71%  '

Set-TextValue $ws.Range('D22') '0.999'
Set-TextValue $ws.Range('E22') '  -0.08%  '

Set-TextValue $ws.Range('E23') '  +10.91%  '

Set-TextValue $ws.Range('D24') '67.65'
Set-TextValue $ws.Range('E24') '  +1.31%  '

Set-TextValue $ws.Range('D25') '630.59'
Set-TextValue $ws.Range('E25') '  +14.08%  '

Set-TextValue $ws.Range('E26') '  +14.03%  '

Set-TextValue $ws.Range('D27') '8.80'
Set-TextValue $ws.Range('E27') '  +0.71%  '

Set-TextValue $ws.Range('D28') '2.599.06'
Set-TextValue $ws.Range('E28') '  +2.38%  '

Set-TextValue $ws.Range('E29') '  +9.91%  '

Set-TextValue $ws.Range('D30') '8.48'
Set-TextValue $ws.Range('E30') '  +3.52%  '

Set-TextValue $ws.Range('E31') '  -0.28%  '

Set-TextValue $ws.Range('E32') '  -1.66%  '

Set-TextValue $ws.Range('E33') '  +2.16%  '

Set-TextValue $ws.Range('D34') '5.23'
Set-TextValue $ws.Range('E34') '  +10.53%  '

Set-TextValue $ws.Range('E35') '  +3.76%  '

Set-TextValue $ws.Range('D36') '0.997'
Set-TextValue $ws.Range('E36') '  -0.26%  '

Set-TextValue $ws.Range('E37') '  +2.26%  '

Set-TextValue $ws.Range('E38') '  +2.99%  '

Set-TextValue $ws.Range('D39') '19.03'
Set-TextValue $ws.Range('E39') '  +2.67%  '

Set-TextValue $ws.Range('E40') '  +3.79%  '

Set-TextValue $ws.Range('D41') '146.76'
Set-TextValue $ws.Range('E41') '  -4.40%  '

Set-TextValue $ws.Range('D42') '2.67'
Set-TextValue $ws.Range('E42') '  +19.71%  '

Set-TextValue $ws.Range('E43') '  -0.02%  '

Set-TextValue $ws.Range('D44') '151.31'
Set-TextValue $ws.Range('E44') '  +2.95%  '

Set-TextValue $ws.Range('E45') '  +4.10%  '

Set-TextValue $ws.Range('D46') '21.28'
Set-TextValue $ws.Range('E46') '  +7.59%  '

Set-TextValue $ws.Range('D47') '0.0552'
Set-TextValue $ws.Range('E47') '  +4.71%  '

Set-TextValue $ws.Range('E48') '  +3.45%  '

Set-TextValue $ws.Range('E49') '  +6.04%  '

Set-TextValue $ws.Range('D50') '0.0930'
Set-TextValue $ws.Range('E50') '  +0.97%  '

Set-TextValue $ws.Range('D51') '0.751'
Set-TextValue $ws.Range('E51') '  +5.75%  '
